# Apply updated crypto price/volume figures to Sheet1 (columns D and E).
# Numeric-looking values in column D are prefixed with a literal apostrophe so
# Excel stores them as text (preserving trailing/leading zeros) instead of
# silently converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.401.32'
$ws.Range("E2").Value = '  +0.08%  '
$ws.Range("D3").Value = '1.847.68'
$ws.Range("D4").Value = '''0.9995'
$ws.Range("E4").Value = '  +0.08%  '
$ws.Range("D5").Value = '''240.25'
$ws.Range("E5").Value = '  -0.19%  '
$ws.Range("D6").Value = '''0.6291'
$ws.Range("E6").Value = '  -1.41%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '''0.07607'
$ws.Range("E8").Value = '  +0.62%  '
$ws.Range("D9").Value = '''0.2931'
$ws.Range("E9").Value = '  -1.24%  '
$ws.Range("E10").Value = '  -1.18%  '
$ws.Range("D11").Value = '''0.07743'
$ws.Range("E11").Value = '  +0.09%  '
$ws.Range("D12").Value = '1.843.47'
$ws.Range("E12").Value = '  -7.13%  '
$ws.Range("E13").Value = '  +0.20%  '
$ws.Range("D14").Value = '''0.00001086'
$ws.Range("E14").Value = '  +9.06%  '
$ws.Range("D15").Value = '''0.6796'
$ws.Range("E15").Value = '  -0.75%  '
$ws.Range("D16").Value = '''83.76'
$ws.Range("D17").Value = '2.102.76'
$ws.Range("E17").Value = '  -7.14%  '
$ws.Range("D18").Value = '''6.182'
$ws.Range("E18").Value = '  +0.02%  '
$ws.Range("D19").Value = '29.417.22'
$ws.Range("E19").Value = '  +0.05%  '
$ws.Range("D20").Value = '''228.92'
$ws.Range("E20").Value = '  -0.24%  '
$ws.Range("E21").Value = '  -0.24%  '
$ws.Range("D23").Value = '''7.464'
$ws.Range("E23").Value = '  -1.53%  '
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("D25").Value = '''157.39'
$ws.Range("E25").Value = '  +0.45%  '
$ws.Range("D26").Value = '''0.1399'
$ws.Range("E26").Value = '  -0.60%  '
$ws.Range("D27").Value = '''8.365'
$ws.Range("E27").Value = '  -0.32%  '
$ws.Range("D29").Value = '''1.465'
$ws.Range("E29").Value = '  -0.29%  '
$ws.Range("D30").Value = '''1.297'
$ws.Range("E30").Value = '  +3.95%  '
$ws.Range("D31").Value = '''0.05580'
$ws.Range("E31").Value = '  -2.21%  '
$ws.Range("E32").Value = '  -0.72%  '
$ws.Range("E33").Value = '  -0.12%  '
$ws.Range("D34").Value = '''1.843'
$ws.Range("E34").Value = '  -0.29%  '
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("D36").Value = '''0.7096'
$ws.Range("E36").Value = '  -1.05%  '
$ws.Range("D37").Value = '''2.584'
$ws.Range("E37").Value = '  -0.26%  '
$ws.Range("D38").Value = '1.233.85'
$ws.Range("E38").Value = '  -1.77%  '
$ws.Range("D39").Value = '''0.01801'
$ws.Range("E39").Value = '  -0.59%  '
$ws.Range("D40").Value = '''2.765'
$ws.Range("E40").Value = '  -0.90%  '
$ws.Range("D41").Value = '''6.440'
$ws.Range("E41").Value = '  +5.70%  '
$ws.Range("D42").Value = '''0.9063'
$ws.Range("E42").Value = '  -0.38%  '
$ws.Range("D43").Value = '''1.0000'
$ws.Range("E44").Value = '  +0.15%  '
$ws.Range("D45").Value = '''66.12'
$ws.Range("E45").Value = '  -0.64%  '
$ws.Range("D46").Value = '''0.00000000121'
$ws.Range("E46").Value = '  +2.34%  '
$ws.Range("D47").Value = '''7.170'
$ws.Range("E47").Value = '  +1.46%  '
$ws.Range("D48").Value = '''0.4022'
$ws.Range("E48").Value = '  -0.25%  '
$ws.Range("D49").Value = '''8.989'
$ws.Range("E49").Value = '  -2.38%  '
$ws.Range("D50").Value = '''1.678'
$ws.Range("E50").Value = '  -1.57%  '
